$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "35.248.41"
$ws.Range("E2").Value = "  +0.26%  "

# Row 3
$ws.Range("D3").Value = "1.895.22"
$ws.Range("E3").Value = "  +2.10%  "

# Row 4
$ws.Range("E4").Value = "  -0.17%  "

# Row 5
$ws.Range("D5").Value = "244.04"
$ws.Range("E5").Value = "  +2.34%  "

# Row 6
$ws.Range("D6").Value = "0.654"
$ws.Range("E6").Value = "  +5.25%  "

# Row 7
$ws.Range("E7").Value = "  -0.37%  "

# Row 8
$ws.Range("D8").Value = "41.40"
$ws.Range("E8").Value = "  -1.58%  "

# Row 9
$ws.Range("D9").Value = "0.351"
$ws.Range("E9").Value = "  +7.40%  "

# Row 10
$ws.Range("D10").Value = "51.97"
$ws.Range("E10").Value = "  +11.60%  "

# Row 12
$ws.Range("E12").Value = "  +0.34%  "

# Row 13
$ws.Range("D13").Value = "2.169.19"
$ws.Range("E13").Value = "  +2.05%  "

# Row 14
$ws.Range("D14").Value = "12.00"
$ws.Range("E14").Value = "  +4.97%  "

# Row 15
$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").Value = "0.693"
$ws.Range("E15").Value = "  +2.50%  "

# Row 16
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "1.888.61"
$ws.Range("E16").Value = "  +1.55%  "

# Row 17
$ws.Range("D17").Value = "4.84"
$ws.Range("E17").Value = "  +2.55%  "

# Row 18
$ws.Range("D18").Value = "35.247.01"
$ws.Range("E18").Value = "  +0.36%  "

# Row 19
$ws.Range("D19").Value = "71.24"

# Row 21
$ws.Range("D21").Value = "239.97"
$ws.Range("E21").Value = "  -0.41%  "

# Row 22
$ws.Range("D22").Value = "12.46"
$ws.Range("E22").Value = "  +2.01%  "

# Row 23
$ws.Range("D23").Value = "4.77"
$ws.Range("E23").Value = "  +0.97%  "

# Row 24
$ws.Range("E24").Value = "  -0.37%  "

# Row 25
$ws.Range("D25").Value = "2.43"
$ws.Range("E25").Value = "  +29.81%  "

# Row 26
$ws.Range("E26").Value = "  +0.59%  "

# Row 27
$ws.Range("E27").Value = "  +1.37%  "

# Row 28
$ws.Range("D28").Value = "8.48"
$ws.Range("E28").Value = "  +6.43%  "

# Row 29
$ws.Range("D29").Value = "18.26"
$ws.Range("E29").Value = "  +3.63%  "

# Row 30
$ws.Range("E30").Value = "  +2.40%  "

# Row 31
$ws.Range("D31").Value = "4.12"
$ws.Range("E31").Value = "  +3.57%  "

# Row 32
$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").Value = "0.0562"
$ws.Range("E32").Value = "  +1.33%  "

# Row 33
$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").Value = "0.936"
$ws.Range("E33").Value = "  +11.35%  "

# Row 34
$ws.Range("E34").Value = "  -0.19%  "

# Row 35
$ws.Range("D35").Value = "4.12"
$ws.Range("E35").Value = "  +3.26%  "

# Row 36
$ws.Range("D36").Value = "1.73"
$ws.Range("E36").Value = "  -3.64%  "

# Row 37
$ws.Range("D37").Value = "2.02"
$ws.Range("E37").Value = "  -0.17%  "

# Row 38
$ws.Range("E38").Value = "  +1.45%  "

# Row 39
$ws.Range("D39").Value = "0.0209"
$ws.Range("E39").Value = "  +4.39%  "

# Row 40
$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").Value = "0.0651"
$ws.Range("E40").Value = "  +17.24%  "

# Row 41
$ws.Range("B41").Value = "ARBITRUM"
$ws.Range("C41").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D41").Value = "1.09"
$ws.Range("E41").Value = "  +2.04%  "

# Row 42
$ws.Range("D42").Value = "16.23"
$ws.Range("E42").Value = "  +9.43%  "

# Row 43
$ws.Range("D43").Value = "89.67"
$ws.Range("E43").Value = "  -0.69%  "

# Row 44
$ws.Range("D44").Value = "1.338.77"
$ws.Range("E44").Value = "  -0.26%  "

# Row 45
$ws.Range("D45").Value = "2.38"
$ws.Range("E45").Value = "  +2.68%  "

# Row 46
$ws.Range("D46").Value = "47.76"
$ws.Range("E46").Value = "  +37.78%  "

# Row 47
$ws.Range("B47").Value = "MXToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D47").Value = "2.78"
$ws.Range("E47").Value = "  +1.73%  "

# Row 48
$ws.Range("B48").Value = "HuobiToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D48").Value = "2.40"
$ws.Range("E48").Value = "  -0.59%  "

# Row 49
$ws.Range("D49").Value = "6.56"
$ws.Range("E49").Value = "  +1.32%  "

# Row 50
$ws.Range("D50").Value = "2.079.67"
$ws.Range("E50").Value = "  +2.07%  "

# Row 51
$ws.Range("D51").Value = "11.14"
$ws.Range("E51").Value = "  -11.13%  "
